# The commit swaps the palette that is actually applied to the deck's
# Slide Master (ppt/theme/theme1.xml): the custom "Integral" colour
# scheme is replaced by the stock Office theme colour scheme (the
# palette that used to sit, unused, behind the Notes Master in
# ppt/theme/theme2.xml).
#
# PowerPoint's object model exposes a theme's 12 colour slots through
# ThemeColorScheme.Colors(i).RGB (1-based: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). COM colour values use the OLE
# 0x00BBGGRR layout, so build each value from its R/G/B bytes instead
# of pasting the "RRGGBB" hex string directly.

function New-RgbValue {
    param(
        [int]$R,
        [int]$G,
        [int]$B
    )
    return $R + ($G * 256) + ($B * 65536)
}

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Target palette: the stock "Office Theme" colour scheme.
$cs.Colors(1).RGB  = New-RgbValue 0x00 0x00 0x00   # dk1      000000
$cs.Colors(2).RGB  = New-RgbValue 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = New-RgbValue 0x44 0x54 0x6A   # dk2      44546A
$cs.Colors(4).RGB  = New-RgbValue 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Colors(5).RGB  = New-RgbValue 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Colors(6).RGB  = New-RgbValue 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Colors(7).RGB  = New-RgbValue 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Colors(8).RGB  = New-RgbValue 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Colors(9).RGB  = New-RgbValue 0x44 0x72 0xC4   # accent5  4472C4
$cs.Colors(10).RGB = New-RgbValue 0x70 0xAD 0x47   # accent6  70AD47
$cs.Colors(11).RGB = New-RgbValue 0x05 0x63 0xC1   # hlink    0563C1
$cs.Colors(12).RGB = New-RgbValue 0x95 0x4F 0x72   # folHlink 954F72
